# PhrasesToRead.xlsx update
# - Adjusts several "#mturk DONE" (column G) counts on Sheet1 (including two
#   rows that previously had no entry at all), which in turn recalculates
#   the "needed?" (column I) formula cells.
# - Clears the AutoFilter's custom filter on column E ("*dour*") while
#   keeping the existing sort state, which also unhides every row that the
#   filter had hidden.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column G updates (existing values changed + two new entries) ---------
$gUpdates = @{
    2   = 0
    4   = 1
    6   = 1
    12  = 0   # previously empty
    19  = 1
    28  = 0
    36  = 1
    39  = 1
    41  = 0
    48  = 1
    49  = 1
    53  = 1
    58  = 0
    70  = 0
    73  = 1
    77  = 1
    83  = 0
    84  = 0
    97  = 0
    99  = 0   # previously empty
    106 = 1
}

foreach ($row in $gUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
}

# --- Remove the autoFilter custom filter on column E (keeps sort state) ---
# This also unhides every row that was hidden because of the filter.
$ws.ShowAllData()
